$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8769625425338745
$ws.Range("B1").Value = 1.448047161102295
$ws.Range("C1").Value = 6.4581298828125
$ws.Range("D1").Value = 1.869495749473572
$ws.Range("E1").Value = 1.08042049407959
